$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all cells we touch so numeric-looking strings
# (e.g. "1.00", "0.0000260", "3.299.88") are preserved exactly as text,
# matching the original inlineStr cell content instead of being
# auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.799.43"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.299.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.38"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.79"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.294.20"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.570"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.572"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.40"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.828.15"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.33"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "575.09"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -8.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.638.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.295.25"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.50"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.75"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.882"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.69"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.94"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.30"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -8.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.90"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.67"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.21"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.57%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.25"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.31"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.56"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.02%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.66"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.86%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "558.47"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.40%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.745.45"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.102"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.47"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.01"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.126"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.10"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.35%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0672"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.54"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.329"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0403"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.03"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -10.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.48"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.57%  "
